$d = $word.ActiveDocument

# --- Change 1: "For example see screenshot below." -> "For example, see screenshot below." ---
# (also removes the now-stale gramStart/gramEnd proofing marks that bracketed "example")
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("For example see screenshot below.", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "For example, see screenshot below.", 2)

# --- Change 2: split " folder. You can check ... directory." so a new clause is inserted ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("folder. You can check this using list (‘ls’) which shows all the files in your current working directory.", `
                              $true, $false, $false, $false, $false, $true, 1, $false, `
                              "folder. If you like, you can check this using list (‘ls’) which shows all the files in your current working directory.", 2)

Write-Output "change1 found=$found1 change2 found=$found2"
